{"js": "// The document has a table describing the \"D\u00f2ng l\u1ec7nh t\u1ea1o m\u1edbi db\" commands.\n// Its header row reads \"Ch\u1ee9c n\u0103ng\" / \"C\u00e2u l\u1ec7nh\" and is followed by a row\n// for \"T\u1ea1o m\u1edbi\". We need to insert a new row, right before \"T\u1ea1o m\u1edbi\"\n// (i.e. right after the header row), for the new \"C\u00e0i \u0111\u1eb7t\" step.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet targetRow = null;\nlet headerRow = null;\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    row.load(\"values\");\n    await context.sync();\n\n    const firstCellText = (row.values[0][0] || \"\").trim();\n    if (firstCellText === \"Ch\u1ee9c n\u0103ng\") {\n      headerRow = row;\n    }\n    if (firstCellText === \"T\u1ea1o m\u1edbi\") {\n      targetRow = row;\n      break;\n    }\n  }\n\n  if (targetRow || headerRow) {\n    break;\n  }\n}\n\nconst newRowValues = [\n  [\"C\u00e0i \u0111\u1eb7t\", \"dotnet tool update --global dotnet-ef --version 7.0.0\"],\n];\n\nif (targetRow) {\n  // Insert the new row directly above the \"T\u1ea1o m\u1edbi\" row.\n  targetRow.insertRows(\"Before\", 1, newRowValues);\n} else if (headerRow) {\n  // Fallback: insert right after the header row.\n  headerRow.insertRows(\"After\", 1, newRowValues);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The table on this page lists the \"D\u00f2ng l\u1ec7nh t\u1ea1o m\u1edbi db\" commands; its\n# header row is \"Ch\u1ee9c n\u0103ng\" / \"C\u00e2u l\u1ec7nh\" and the first data row is \"T\u1ea1o m\u1edbi\".\n# We need to insert a new row right after the header (i.e. right before the\n# \"T\u1ea1o m\u1edbi\" row) holding the new \"C\u00e0i \u0111\u1eb7t\" step.\n\n$targetTable = $null\n$targetRowIndex = 0\n$headerRowIndex = 0\n\nfor ($t = 1; $t -le $d.Tables.Count; $t++) {\n    $table = $d.Tables.Item($t)\n    for ($i = 1; $i -le $table.Rows.Count; $i++) {\n        $cellText = $table.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13)\n        if ($cellText -eq \"Ch\u1ee9c n\u0103ng\") {\n            $headerRowIndex = $i\n        }\n        if ($cellText -eq \"T\u1ea1o m\u1edbi\") {\n            $targetTable = $table\n            $targetRowIndex = $i\n            break\n        }\n    }\n    if ($targetRowIndex -gt 0) {\n        break\n    }\n    if ($headerRowIndex -gt 0) {\n        $targetTable = $table\n        break\n    }\n}\n\nif ($targetTable -ne $null) {\n    if ($targetRowIndex -gt 0) {\n        # Insert the new row directly above the \"T\u1ea1o m\u1edbi\" row.\n        $newRow = $targetTable.Rows.Add($targetTable.Rows.Item($targetRowIndex))\n    } else {\n        # Fallback: insert right after the header row.\n        $after = $headerRowIndex + 1\n        if ($after -gt $targetTable.Rows.Count) {\n            $newRow = $targetTable.Rows.Add()\n        } else {\n            $newRow = $targetTable.Rows.Add($targetTable.Rows.Item($after))\n        }\n    }\n\n    $newRow.Cells.Item(1).Range.Text = \"C\u00e0i \u0111\u1eb7t\"\n    $newRow.Cells.Item(2).Range.Text = \"dotnet tool update --global dotnet-ef --version 7.0.0\"\n}\n"}
